$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 10: Inscritos (E10) 26 -> 27
$ws.Range("E10").Value = 27

# Row 17: Inscritos (E17) 19 -> 20, Pagos (F17) 11 -> 12, Inscricoes homologadas (H17) 11 -> 12
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = 12
